# Update EPEX Spot prices workbook with the latest day's data
# (27-aug on "Prix Spot", 2025-08-25 on "Gaz" and "CO2")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Prix Spot": append a new date column (BW) after BV
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting (bold header, borders, alignment) of the last
# existing column (BV) into the new column (BW) before writing values.
$ws1.Range("BV1:BV25").Copy($ws1.Range("BW1:BW25")) | Out-Null

$ws1.Range("BW1").Value = "27-aug"

$pricesSpot = @{
    2  = 97.41
    3  = 90.02
    4  = 86.55
    5  = 80.36
    6  = 79.04000000000001
    7  = 85.90000000000001
    8  = 102.84
    9  = 113.12
    10 = 113.64
    11 = 101.63
    12 = 85.20999999999999
    13 = 77.75
    14 = 63.65
    15 = 49.85
    16 = 50.62
    17 = 50.62
    18 = 63.81
    19 = 78.73
    20 = 89.61
    21 = 99
    22 = 109.1
    23 = 114.32
    24 = 113.63
    25 = 103.57
}

foreach ($row in $pricesSpot.Keys) {
    $ws1.Cells.Item($row, 75).Value = $pricesSpot[$row]
}

# ---------------------------------------------------------------
# Sheet "Gaz": append a new row (72) with the latest closing price
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")

# Force text formatting on the date cell first so Excel doesn't
# auto-convert the "yyyy-mm-dd" string into a date serial number,
# then restore the default "Normal" style so no style index is left
# attached to the cell (matches the other date cells in the column).
$ws2.Range("A72").NumberFormat = "@"
$ws2.Range("A72").Value = "2025-08-25"
$ws2.Range("A72").Style = "Normal"
$ws2.Range("B72").Value = 32.625

# ---------------------------------------------------------------
# Sheet "CO2": append a new row (72) with the latest closing price
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A72").NumberFormat = "@"
$ws3.Range("A72").Value = "2025-08-25"
$ws3.Range("A72").Style = "Normal"
$ws3.Range("B72").Value = 71.52
